$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "10 to 64 feet" header to the correct "40 to 64 feet"
$ws.Range("D1").Value = "40 to 64 feet"

# Update the active selection on the sheet (as saved in the view state)
$ws.Range("G1").Select()
